# Sari_R3_3 template update
#  - Week-number column (B7:B59) now filled with the epidemiological week index (1-53)
#  - Header cell C6 gets the same wrap-text style as the other header cells
#  - Selection moved to the top title row (B1:O1) as it was left by the author
#  - Value axis of the stacked bar chart gets a title ("Número de casos SARI")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet1 data: week numbers 1..53 in column B (rows 7-59) ---------------
for ($i = 1; $i -le 53; $i++) {
    $ws.Cells.Item($i + 6, 2).Value = $i
}

# --- Header row: C6 should wrap text like the neighboring headers ----------
$ws.Range("C6").WrapText = $true

# --- Selection left on the title band after editing ------------------------
$ws.Range("B1:O1").Select() | Out-Null

# --- Chart: add a title to the value (Y) axis -------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$valueAxis = $chart.Axes(2)          # 2 = xlValue
$valueAxis.HasTitle = $true
$valueAxis.AxisTitle.Text = "Número de casos SARI"
